# Apply styles.xml changes described by the commit:
#   "Meg and I added content on the 'for educators' page"
#
# 1. New "Abstract Title" paragraph style, inserted just before "Abstract".
# 2. "Abstract" style: spacing-before 300 -> 100 twips.
# 3. New "Footnote Block Text" paragraph style, based on "Footnote Text".
# 4. Pandoc syntax-highlighting character styles recolored/bolded:
#    ConstantTok, SpecialCharTok, FunctionTok, AttributeTok.

$d = $word.ActiveDocument

# --- 1. Abstract Title (new paragraph style) ---------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.SpaceBefore = 15

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060   # RGB(0x345A8A)

# --- 2. Abstract: spacing before 300 -> 100 -----------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. Footnote Block Text (new paragraph style) -----------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0

# --- 4. Pandoc syntax-highlight token colors -----------------------------
$constantTok = $d.Styles("ConstantTok")
$constantTok.Font.Color = 153999      # RGB(0x8F5902)

$specialCharTok = $d.Styles("SpecialCharTok")
$specialCharTok.Font.Color = 23758    # RGB(0xCE5C00)
$specialCharTok.Font.Bold = $true

$functionTok = $d.Styles("FunctionTok")
$functionTok.Font.Color = 8866336     # RGB(0x204A87)
$functionTok.Font.Bold = $true

$attributeTok = $d.Styles("AttributeTok")
$attributeTok.Font.Color = 8866336    # RGB(0x204A87)

Write-Output "styles updated"
